# Updated cryptos list with GitHub Actions
# Applies the Price (column D) and Volume(1h) (column E) updates described
# by the diff to the active worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.993.21"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.681.72"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "215.49"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "0.516"
$ws.Range("E6").Value = "  -2.53%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "0.0619"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "21.02"
$ws.Range("E10").Value = "  +3.20%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.919.19"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "1.698.28"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("E15").Value = "  +1.63%  "
$ws.Range("D16").Value = "65.94"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "27.018.18"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "8.16"
$ws.Range("E18").Value = "  +4.30%  "
$ws.Range("D19").Value = "236.34"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").Value = "2.12"
$ws.Range("E24").Value = "  -4.32%  "
$ws.Range("D25").Value = "146.46"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -3.25%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Value = "1.497.20"
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").Value = "1.68"
$ws.Range("E35").Value = "  +3.92%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "0.585"
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").Value = "  +3.32%  "
$ws.Range("E40").Value = "  +7.32%  "
$ws.Range("E41").Value = "  -4.65%  "
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").Value = "67.53"
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "1.824.00"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").Value = "90.32"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("E49").Value = "  +3.42%  "
$ws.Range("D50").Value = "7.80"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("E51").Value = "  -0.01%  "
